$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B first (row order)
$ws.Range("B9").Value = "/statistique/"
$ws.Range("B10").Value = 1
$ws.Range("B11").Value = 2
$ws.Range("B12").Value = 3
$ws.Range("B13").Value = 4
$ws.Range("B14").Value = 5
$ws.Range("B15").Value = 6

# Column C next (row order) - introduces new shared strings 19-24
$ws.Range("C10").Value = "le chiffre affaire total"
$ws.Range("C11").Value = "le benefice total (affaire - salaire)"
$ws.Range("C12").Value = "le chiffre affaire max (nom de l 'employee + montant)"
$ws.Range("C13").Value = "le chiffre affaire  min  (nom de l 'employee + montant)"
$ws.Range("C14").Value = "employee le plus rentable"
$ws.Range("C15").Value = "employee les moins rentable (liste <nom de l 'employee + montant> des employe dont le benefice < benefice moyen )"

# Column D last (row order) - introduces new shared strings 25-27
$ws.Range("D10").Value = "OK, à verifier quantité"
$ws.Range("D11").Value = "OK, à verifier quantité"
$ws.Range("D12").Value = "OK, je retourne un obj Employee, entier!"
$ws.Range("D13").Value = "OK, je retourne un obj Employee, entier!"
$ws.Range("D14").Value = "KO"
$ws.Range("D15").Value = "KO"

# Highlight D14:D15 in red
$ws.Range("D14").Interior.Color = 255
$ws.Range("D15").Interior.Color = 255

# Column C width change (target stored width 108.140625; engine quantizes
# ColumnWidth to 1/6-character steps, so 107.33 is the closest settable
# value that rounds to the nearest representable stored width, 108.1667)
$ws.Columns("C").ColumnWidth = 107.33

# Selection change
$ws.Range("C17").Select()
